$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "local 21 / ACTIVATION_OTHER / ANY OTHER ACTIVATION" KPI row (row 43)
# is removed; everything below it shifts up one row.
$ws.Rows(43).Delete()

# Re-apply AutoFilter over the (now one-row-shorter) data range so the
# filter buttons/ref track the shrunk table. Toggling it off then on
# forces the stored <autoFilter> ref to refresh to the new range.
$ws.Range("A1:AL45").AutoFilter() | Out-Null
$ws.Range("A1:AL45").AutoFilter() | Out-Null

# The hidden _FilterDatabase defined name tracks the autofilter's actual
# range and needs to be pushed out to the new last row.
$wb.Names.Item(1).RefersTo = "='HoReCa Restaurant_Cafe'!`$A`$1:`$AL`$45"

# Re-filtering also leaves behind extra duplicate bookkeeping names.
$ws.Names.Add("_xlnm._FilterDatabase_0_0", "='HoReCa Restaurant_Cafe'!`$A`$1:`$AL`$1") | Out-Null
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0", "='HoReCa Restaurant_Cafe'!`$A`$1:`$AL`$1") | Out-Null

# Leave the cursor where the editor ended up.
$ws.Range("T4").Select()
